$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.171875
$ws.Range("C2").Value = 0.615625
$ws.Range("P2").Value = 0.121875
$ws.Range("S2").Value = 0.090625
$ws.Range("B3").Value = 0.004901960784313725
$ws.Range("C3").Value = 0.02450980392156863
$ws.Range("J3").Value = 0.0196078431372549
$ws.Range("P3").Value = 0.7696078431372549
$ws.Range("S3").Value = 0.1813725490196078
$ws.Range("O4").Value = 0.025
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.225
$ws.Range("B6").Value = 0.0380952380952381
$ws.Range("D6").Value = 0.004761904761904762
$ws.Range("F6").Value = 0.06666666666666667
$ws.Range("J6").Value = 0.2428571428571429
$ws.Range("O6").Value = 0.01428571428571429
$ws.Range("Q6").Value = 0.1333333333333333
$ws.Range("R6").Value = 0.1238095238095238
$ws.Range("S6").Value = 0.3761904761904762
$ws.Range("B7").Value = 0.1139896373056995
$ws.Range("D7").Value = 0.005181347150259068
$ws.Range("F7").Value = 0.04145077720207254
$ws.Range("J7").Value = 0.09326424870466321
$ws.Range("O7").Value = 0.02590673575129534
$ws.Range("Q7").Value = 0.1191709844559585
$ws.Range("R7").Value = 0.1139896373056995
$ws.Range("S7").Value = 0.4870466321243523
$ws.Range("B8").Value = 0.1334661354581673
$ws.Range("D8").Value = 0.01593625498007968
$ws.Range("E8").Value = 0.00597609561752988
$ws.Range("F8").Value = 0.049800796812749
$ws.Range("J8").Value = 0.08964143426294821
$ws.Range("O8").Value = 0.02390438247011952
$ws.Range("Q8").Value = 0.151394422310757
$ws.Range("R8").Value = 0.1175298804780877
$ws.Range("S8").Value = 0.4123505976095618
$ws.Range("B9").Value = 0.09644670050761421
$ws.Range("D9").Value = 0.03553299492385787
$ws.Range("E9").Value = 0.005076142131979695
$ws.Range("F9").Value = 0.03553299492385787
$ws.Range("J9").Value = 0.07614213197969544
$ws.Range("O9").Value = 0.02030456852791878
$ws.Range("Q9").Value = 0.1065989847715736
$ws.Range("R9").Value = 0.09644670050761421
$ws.Range("S9").Value = 0.5279187817258884
$ws.Range("B10").Value = 0.1296296296296296
$ws.Range("D10").Value = 0.02037037037037037
$ws.Range("E10").Value = 0.002777777777777778
$ws.Range("F10").Value = 0.06759259259259259
$ws.Range("J10").Value = 0.08796296296296297
$ws.Range("O10").Value = 0.02777777777777778
$ws.Range("Q10").Value = 0.1490740740740741
$ws.Range("R10").Value = 0.1185185185185185
$ws.Range("S10").Value = 0.3962962962962963
$ws.Range("G11").Value = 0.1111111111111111
$ws.Range("J11").Value = 0.08602150537634409
$ws.Range("K11").Value = 0.1541218637992832
$ws.Range("L11").Value = 0.6200716845878136
$ws.Range("S11").Value = 0.02867383512544803
$ws.Range("G12").Value = 0.6885245901639344
$ws.Range("J12").Value = 0.2131147540983606
$ws.Range("L12").Value = 0.01639344262295082
$ws.Range("S12").Value = 0.08196721311475409
$ws.Range("F15").Value = 0.02325581395348837
$ws.Range("H15").Value = 0.1767441860465116
$ws.Range("I15").Value = 0.06976744186046512
$ws.Range("J15").Value = 0.3209302325581395
$ws.Range("K15").Value = 0.03720930232558139
$ws.Range("M15").Value = 0.0186046511627907
$ws.Range("O15").Value = 0.04186046511627907
$ws.Range("S15").Value = 0.3116279069767442
$ws.Range("F16").Value = 0.02252252252252252
$ws.Range("H16").Value = 0.1486486486486487
$ws.Range("I16").Value = 0.1081081081081081
$ws.Range("J16").Value = 0.3513513513513514
$ws.Range("K16").Value = 0.06756756756756757
$ws.Range("M16").Value = 0.01801801801801802
$ws.Range("N16").Value = 0.009009009009009009
$ws.Range("O16").Value = 0.06756756756756757
$ws.Range("S16").Value = 0.2072072072072072
$ws.Range("F17").Value = 0.02622950819672131
$ws.Range("H17").Value = 0.1967213114754098
$ws.Range("I17").Value = 0.09508196721311475
$ws.Range("J17").Value = 0.3573770491803279
$ws.Range("K17").Value = 0.09836065573770492
$ws.Range("M17").Value = 0.02295081967213115
$ws.Range("O17").Value = 0.04918032786885246
$ws.Range("S17").Value = 0.1540983606557377
$ws.Range("F18").Value = 0.04280155642023346
$ws.Range("H18").Value = 0.1828793774319066
$ws.Range("I18").Value = 0.08560311284046693
$ws.Range("J18").Value = 0.3618677042801556
$ws.Range("K18").Value = 0.11284046692607
$ws.Range("M18").Value = 0.01945525291828794
$ws.Range("O18").Value = 0.05447470817120623
$ws.Range("S18").Value = 0.1400778210116732
$ws.Range("F19").Value = 0.02121433796634967
$ws.Range("H19").Value = 0.2392099487929773
$ws.Range("I19").Value = 0.07973664959765911
$ws.Range("J19").Value = 0.324067300658376
$ws.Range("K19").Value = 0.1068032187271397
$ws.Range("M19").Value = 0.02852962692026335
$ws.Range("O19").Value = 0.0592538405267008
$ws.Range("S19").Value = 0.141185076810534
